$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised historical values (rows 226-251) ---
$ws.Range("C226:F226").Value = 1486571000000
$ws.Range("C230:F230").Value = 1496040000000
$ws.Range("C237:F237").Value = 1485900000000
$ws.Range("C238:F238").Value = 1498300000000
$ws.Range("C239:F239").Value = 1517100000000
$ws.Range("C240:F240").Value = 1563100000000
$ws.Range("C241:F241").Value = 1563400000000
$ws.Range("C242:F242").Value = 1570300000000
$ws.Range("C243:F243").Value = 1593400000000
$ws.Range("C244:F244").Value = 1567300000000
$ws.Range("C245:F245").Value = 1568100000000
$ws.Range("C246:F246").Value = 1622200000000
$ws.Range("C247:F247").Value = 1606000000000
$ws.Range("C248:F248").Value = 1627700000000
$ws.Range("C249:F249").Value = 1645600000000
$ws.Range("C250:F250").Value = 1629400000000
$ws.Range("C251:F251").Value = 1678100000000

# --- Append new rows 258-260 ---
$ws.Range("A257:G257").Copy()
$ws.Range("A258:G258").PasteSpecial(-4122)
$ws.Range("A258").Value = 45078.41666666666
$ws.Range("B258").Value = "ECONOMICS:AEM2"
$ws.Range("C258:F258").Value = 1855306000000
$ws.Range("G258").Value = 0

$ws.Range("A258:G258").Copy()
$ws.Range("A259:G259").PasteSpecial(-4122)
$ws.Range("A259").Value = 45108.41666666666
$ws.Range("B259").Value = "ECONOMICS:AEM2"
$ws.Range("C259:F259").Value = 1858844000000
$ws.Range("G259").Value = 0

$ws.Range("A259:G259").Copy()
$ws.Range("A260:G260").PasteSpecial(-4122)
$ws.Range("A260").Value = 45139.41666666666
$ws.Range("B260").Value = "ECONOMICS:AEM2"
$ws.Range("C260:F260").Value = 1860330000000
$ws.Range("G260").Value = 0

$excel.CutCopyMode = 0
